$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (old D:K shift right to F:M).
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy the number/date formatting from the (now-shifted) old column D/E
# range (F:G) onto the freshly inserted D:E columns so the new cells pick
# up the same styles (date format / thousands format) as their neighbours.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 4147000
$ws.Range("E8").Value = 4888000
$ws.Range("D9").Value = 2387000
$ws.Range("E9").Value = 2820000
$ws.Range("D10").Value = 1760000
$ws.Range("E10").Value = 2068000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 10000
$ws.Range("E14").Value = 27000
$ws.Range("D15").Value = 57000
$ws.Range("E15").Value = 57000
$ws.Range("D17").Value = 3514000
$ws.Range("E17").Value = 4101000
$ws.Range("D18").Value = 633000
$ws.Range("E18").Value = 787000
$ws.Range("D20").Value = 6000
$ws.Range("E20").Value = 5000
$ws.Range("D21").Value = 841000
$ws.Range("E21").Value = 993000
$ws.Range("D22").Value = 48000
$ws.Range("E22").Value = 54000
$ws.Range("D23").Value = 591000
$ws.Range("E23").Value = 738000
$ws.Range("D24").Value = 224000
$ws.Range("E24").Value = 262000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 367000
$ws.Range("E26").Value = 476000
$ws.Range("D27").Value = 365000
$ws.Range("E27").Value = 471000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 100000
$ws.Range("E29").Value = 146000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -6000
$ws.Range("E32").Value = -5000
$ws.Range("D33").Value = 465000
$ws.Range("E33").Value = 617000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 465000
$ws.Range("E35").Value = 617000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1248000
$ws.Range("E41").Value = 1093000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2733000
$ws.Range("E43").Value = 3023000
$ws.Range("D44").Value = 1980000
$ws.Range("E44").Value = 1813000
$ws.Range("D45").Value = 697000
$ws.Range("E45").Value = 690000
$ws.Range("D46").Value = 6658000
$ws.Range("E46").Value = 6619000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 3551000
$ws.Range("E48").Value = 3562000
$ws.Range("D49").Value = 9182000
$ws.Range("E49").Value = 11957000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1038000
$ws.Range("E52").Value = 1003000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 20429000
$ws.Range("E54").Value = 20390000
$ws.Range("D57").Value = 1794000
$ws.Range("E57").Value = 1943000
$ws.Range("D58").Value = 3320000
$ws.Range("E58").Value = 1623000
$ws.Range("D59").Value = 2426000
$ws.Range("E59").Value = 2598000
$ws.Range("D60").Value = 7540000
$ws.Range("E60").Value = 6164000
$ws.Range("D61").Value = 2641000
$ws.Range("E61").Value = 3137000
$ws.Range("D62").Value = 1972000
$ws.Range("E62").Value = 2099000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 12193000
$ws.Range("E66").Value = 11443000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 23252000
$ws.Range("E72").Value = 23072000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 8236000
$ws.Range("E76").Value = 8947000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 465000
$ws.Range("E81").Value = 617000
$ws.Range("D83").Value = 202000
$ws.Range("E83").Value = 201000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 323000
$ws.Range("E89").Value = 1024000
$ws.Range("D91").Value = -155000
$ws.Range("E91").Value = -303000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -259000
$ws.Range("E94").Value = -1788000
$ws.Range("D96").Value = -305000
$ws.Range("E96").Value = -305000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 98000
$ws.Range("E100").Value = -1532000
$ws.Range("D101").Value = -7000
$ws.Range("E101").Value = -22000
$ws.Range("D102").Value = 155000
$ws.Range("E102").Value = -2318000
